$wb = $excel.ActiveWorkbook

# --- 1. "compounds" sheet: selection moves from G28 to B1:G1 (tabSelected
#        flag follows the active sheet automatically, so do this before we
#        add/activate the new sheet below). ---
$wsCompounds = $wb.Worksheets.Item("compounds")
$wsCompounds.Range("B1:G1").Select() | Out-Null

# --- 2. Add the new "components" worksheet after "compounds" (becomes the
#        5th / last tab and the active sheet). ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "components"

# Header row (reuses the same "as-is / almost as-is / implement / not
# supported / ionize / done" strings already used on the other method
# tables).
$ws.Range("B1").Value = "as-is"
$ws.Range("C1").Value = "almost as-is"
$ws.Range("D1").Value = "implement"
$ws.Range("E1").Value = "not supported"
$ws.Range("F1").Value = "ionize"
$ws.Range("G1").Value = "done"

# Method rows.
$ws.Range("A2").Value = "$"
$ws.Range("B2").Value = "X"
$ws.Range("G2").Value = "X"

$ws.Range("A3").Value = "["
$ws.Range("C3").Value = "X"
$ws.Range("G3").Value = "X"

$ws.Range("A4").Value = "[["
$ws.Range("B4").Value = "X"
$ws.Range("G4").Value = "X"

$ws.Range("A5").Value = "as.data.table"
$ws.Range("B5").Value = "X"
$ws.Range("G5").Value = "X"

$ws.Range("A6").Value = "componentInfo"
$ws.Range("B6").Value = "X"
$ws.Range("G6").Value = "X"

$ws.Range("A7").Value = "componentTable"
$ws.Range("B7").Value = "X"
$ws.Range("G7").Value = "X"

$ws.Range("A8").Value = "consensus"
$ws.Range("E8").Value = "X"

$ws.Range("A9").Value = "filter"
$ws.Range("C9").Value = "X"
$ws.Range("G9").Value = "X"

$ws.Range("A10").Value = "findFGroup"
$ws.Range("B10").Value = "X"
$ws.Range("G10").Value = "X"

$ws.Range("A11").Value = "groupNames"
$ws.Range("B11").Value = "X"
$ws.Range("G11").Value = "X"

$ws.Range("A12").Value = "initialize"
$ws.Range("C12").Value = "X"
$ws.Range("G12").Value = "X"

$ws.Range("A13").Value = "length"
$ws.Range("B13").Value = "X"
$ws.Range("G13").Value = "X"

$ws.Range("A14").Value = "names"
$ws.Range("B14").Value = "X"
$ws.Range("G14").Value = "X"

$ws.Range("A15").Value = "plotEIC"
$ws.Range("B15").Value = "X"
$ws.Range("D15").Value = "X"
$ws.Range("G15").Value = "X"
$ws.Range("H15").Value = "Seems enough, assuming we're not planning to merge components"

$ws.Range("A16").Value = "plotEICHash"
$ws.Range("B16").Value = "X"
$ws.Range("G16").Value = "X"

$ws.Range("A17").Value = "plotSpec"
$ws.Range("B17").Value = "X"
$ws.Range("D17").Value = "X"
$ws.Range("G17").Value = "X"
$ws.Range("H17").Value = "Seems enough, assuming we're not planning to merge components"

$ws.Range("A18").Value = "plotSpecHash"
$ws.Range("B18").Value = "X"
$ws.Range("G18").Value = "X"

$ws.Range("A19").Value = "show"
$ws.Range("C19").Value = "X"
$ws.Range("G19").Value = "X"

# Column A width, auto-fit like the other method-table sheets.
$ws.Columns.Item(1).AutoFit()

# Final selection/active cell on the new sheet.
$ws.Range("H17").Select() | Out-Null
